# Auto-generated Excel COM-interop script
# Applies numeric cell value updates to multiple worksheets as described in the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3731200
$ws.Range("J17").Value = 3731200
$ws.Range("L17").Value = 11193600
$ws.Range("N17").Value = -11193936
$ws.Range("H32").Value = 2190.4285
$ws.Range("I32").Value = 1620.2
$ws.Range("J32").Value = 2368.625
$ws.Range("K32").Value = 1620.2
$ws.Range("L32").Value = 2368.625
$ws.Range("M32").Value = -1294.2
$ws.Range("N32").Value = -3020.625
$ws.Range("H43").Value = 2164.5454
$ws.Range("I43").Value = 1850
$ws.Range("J43").Value = 2234.4443
$ws.Range("K43").Value = 1850
$ws.Range("L43").Value = 2234.4443
$ws.Range("M43").Value = -1781
$ws.Range("N43").Value = -2372.4443
$ws.Range("H127").Value = 1605.7142
$ws.Range("I127").Value = 368
$ws.Range("K127").Value = 1104
$ws.Range("M127").Value = 3856
$ws.Range("H134").Value = 26528
$ws.Range("J134").Value = 26528
$ws.Range("L134").Value = 26528
$ws.Range("N134").Value = -36668
$ws.Range("H138").Value = 4564.1274
$ws.Range("I138").Value = 1031.0769
$ws.Range("J138").Value = 7069.382
$ws.Range("K138").Value = 3093.2307
$ws.Range("L138").Value = 21208.146
$ws.Range("M138").Value = 2046.7693
$ws.Range("N138").Value = -31488.146

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1720.1111
$ws.Range("I2").Value = 1124.9286
$ws.Range("J2").Value = 3803.25
$ws.Range("K2").Value = 1124.9286
$ws.Range("L2").Value = 3803.25
$ws.Range("M2").Value = -1011.9286
$ws.Range("N2").Value = -4029.25
$ws.Range("H32").Value = 5141.8955
$ws.Range("I32").Value = 3352.0186
$ws.Range("J32").Value = 12576.77
$ws.Range("K32").Value = 3352.0186
$ws.Range("L32").Value = 12576.77
$ws.Range("M32").Value = -3065.0186
$ws.Range("N32").Value = -13150.77
$ws.Range("H45").Value = 5024.7915
$ws.Range("I45").Value = 6077.2104
$ws.Range("J45").Value = 1025.6
$ws.Range("K45").Value = 6077.2104
$ws.Range("L45").Value = 1025.6
$ws.Range("M45").Value = -5700.2104
$ws.Range("N45").Value = -1779.6
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
$ws.Range("H74").Value = 1336.0233
$ws.Range("I74").Value = 1234
$ws.Range("J74").Value = 1571.4615
$ws.Range("K74").Value = 1234
$ws.Range("L74").Value = 1571.4615
$ws.Range("M74").Value = -360
$ws.Range("N74").Value = -3319.4615
$ws.Range("H77").Value = 1336.0233
$ws.Range("I77").Value = 1234
$ws.Range("J77").Value = 1571.4615
$ws.Range("K77").Value = 6170
$ws.Range("L77").Value = 7857.307499999999
$ws.Range("M77").Value = -1802
$ws.Range("N77").Value = -16593.3075
$ws.Range("H110").Value = 931.0909
$ws.Range("I110").Value = 884.2
$ws.Range("J110").Value = 1400
$ws.Range("K110").Value = 884.2
$ws.Range("L110").Value = 1400
$ws.Range("M110").Value = 1160.8
$ws.Range("N110").Value = -5490
$ws.Range("H116").Value = 1720.1111
$ws.Range("I116").Value = 1124.9286
$ws.Range("J116").Value = 3803.25
$ws.Range("K116").Value = 1124.9286
$ws.Range("L116").Value = 3803.25
$ws.Range("M116").Value = 1169.0714
$ws.Range("N116").Value = -8391.25
$ws.Range("H132").Value = 2583.6365
$ws.Range("J132").Value = 6688
$ws.Range("L132").Value = 20064
$ws.Range("N132").Value = -25124

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1720.1111
$ws.Range("I3").Value = 1124.9286
$ws.Range("J3").Value = 3803.25
$ws.Range("K3").Value = 1124.9286
$ws.Range("L3").Value = 3803.25
$ws.Range("M3").Value = -1010.9286
$ws.Range("N3").Value = -4031.25
$ws.Range("H105").Value = 1539
$ws.Range("I105").Value = 1020.75
$ws.Range("J105").Value = 3266.5
$ws.Range("K105").Value = 1020.75
$ws.Range("L105").Value = 3266.5
$ws.Range("M105").Value = 726.25
$ws.Range("N105").Value = -6760.5
$ws.Range("H134").Value = 3584.577
$ws.Range("I134").Value = 4105.5405
$ws.Range("K134").Value = 12316.6215
$ws.Range("M134").Value = -9781.621500000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1545.8334
$ws.Range("I16").Value = 1140
$ws.Range("J16").Value = 1835.7142
$ws.Range("K16").Value = 1140
$ws.Range("L16").Value = 1835.7142
$ws.Range("M16").Value = -853
$ws.Range("N16").Value = -2409.7142
$ws.Range("H31").Value = 3024.9219
$ws.Range("I31").Value = 1463.8684
$ws.Range("J31").Value = 5306.4614
$ws.Range("K31").Value = 1463.8684
$ws.Range("L31").Value = 5306.4614
$ws.Range("M31").Value = -1168.8684
$ws.Range("N31").Value = -5896.4614
$ws.Range("H34").Value = 3024.9219
$ws.Range("I34").Value = 1463.8684
$ws.Range("J34").Value = 5306.4614
$ws.Range("K34").Value = 1463.8684
$ws.Range("L34").Value = 5306.4614
$ws.Range("M34").Value = -1261.8684
$ws.Range("N34").Value = -5710.4614
$ws.Range("H58").Value = 1167.7966
$ws.Range("I58").Value = 651.85
$ws.Range("J58").Value = 2254
$ws.Range("K58").Value = 651.85
$ws.Range("L58").Value = 2254
$ws.Range("M58").Value = -448.85
$ws.Range("N58").Value = -2660
$ws.Range("H94").Value = 4522.8096
$ws.Range("I94").Value = 4932
$ws.Range("J94").Value = 4271
$ws.Range("K94").Value = 4932
$ws.Range("L94").Value = 4271
$ws.Range("M94").Value = -4481
$ws.Range("N94").Value = -5173
$ws.Range("H105").Value = 2291.9092
$ws.Range("I105").Value = 3000
$ws.Range("J105").Value = 1701.8334
$ws.Range("K105").Value = 3000
$ws.Range("L105").Value = 1701.8334
$ws.Range("M105").Value = -1253
$ws.Range("N105").Value = -5195.8334
$ws.Range("H113").Value = 1545.8334
$ws.Range("I113").Value = 1140
$ws.Range("J113").Value = 1835.7142
$ws.Range("K113").Value = 1140
$ws.Range("L113").Value = 1835.7142
$ws.Range("M113").Value = 1030
$ws.Range("N113").Value = -6175.7142
$ws.Range("H136").Value = 1167.7966
$ws.Range("I136").Value = 651.85
$ws.Range("J136").Value = 2254
$ws.Range("K136").Value = 1955.55
$ws.Range("L136").Value = 6762
$ws.Range("M136").Value = 594.4499999999998
$ws.Range("N136").Value = -11862

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 3846261.2
$ws.Range("I12").Value = 10000105
$ws.Range("J12").Value = 109.125
$ws.Range("K12").Value = 30000315
$ws.Range("L12").Value = 327.375
$ws.Range("M12").Value = -30000142
$ws.Range("N12").Value = -673.375
$ws.Range("H94").Value = 2964.8
$ws.Range("I94").Value = 2024
$ws.Range("J94").Value = 3200
$ws.Range("K94").Value = 6072
$ws.Range("L94").Value = 9600
$ws.Range("M94").Value = -5396
$ws.Range("N94").Value = -10952

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5913.1816
$ws.Range("I70").Value = 5633.3335
$ws.Range("K70").Value = 5633.3335
$ws.Range("M70").Value = -5363.3335
$ws.Range("H73").Value = 5913.1816
$ws.Range("I73").Value = 5633.3335
$ws.Range("K73").Value = 5633.3335
$ws.Range("M73").Value = -4697.3335
$ws.Range("H113").Value = 47620216
$ws.Range("I113").Value = 76923976
$ws.Range("J113").Value = 1600
$ws.Range("K113").Value = 76923976
$ws.Range("L113").Value = 1600
$ws.Range("M113").Value = -76921806
$ws.Range("N113").Value = -5940
$ws.Range("H122").Value = 4323743
$ws.Range("I122").Value = 7203859
$ws.Range("J122").Value = 3569
$ws.Range("K122").Value = 21611577
$ws.Range("L122").Value = 10707
$ws.Range("M122").Value = -21609127
$ws.Range("N122").Value = -15607
$ws.Range("H126").Value = 5536.7417
$ws.Range("I126").Value = 10118.667
$ws.Range("K126").Value = 30356.001
$ws.Range("M126").Value = -27886.001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1881.1428
$ws.Range("I7").Value = 1650.25
$ws.Range("J7").Value = 2620
$ws.Range("K7").Value = 1650.25
$ws.Range("L7").Value = 2620
$ws.Range("M7").Value = -1538.25
$ws.Range("N7").Value = -2844
$ws.Range("H16").Value = 759.8570999999999
$ws.Range("I16").Value = 663.8
$ws.Range("J16").Value = 1000
$ws.Range("K16").Value = 663.8
$ws.Range("L16").Value = 1000
$ws.Range("M16").Value = -493.8
$ws.Range("N16").Value = -1340
$ws.Range("H40").Value = 500002500
$ws.Range("I40").Value = 1000000000
$ws.Range("J40").Value = 4995
$ws.Range("K40").Value = 1000000000
$ws.Range("L40").Value = 4995
$ws.Range("M40").Value = -999999864
$ws.Range("N40").Value = -5267
$ws.Range("H68").Value = 500001250
$ws.Range("I68").Value = 2500
$ws.Range("K68").Value = 2500
$ws.Range("M68").Value = -1751
$ws.Range("H71").Value = 500001250
$ws.Range("I71").Value = 2500
$ws.Range("K71").Value = 12500
$ws.Range("M71").Value = -8756
$ws.Range("H93").Value = 20000922
$ws.Range("I93").Value = 812.7368
$ws.Range("J93").Value = 83334600
$ws.Range("K93").Value = 812.7368
$ws.Range("L93").Value = 83334600
$ws.Range("M93").Value = 435.2632
$ws.Range("N93").Value = -83337096
$ws.Range("H122").Value = 4073838.2
$ws.Range("I122").Value = 6495516.5
$ws.Range("J122").Value = 1114009.4
$ws.Range("K122").Value = 19486549.5
$ws.Range("L122").Value = 3342028.2
$ws.Range("M122").Value = -19484099.5
$ws.Range("N122").Value = -3346928.2
$ws.Range("H126").Value = 1881.1428
$ws.Range("I126").Value = 1650.25
$ws.Range("J126").Value = 2620
$ws.Range("K126").Value = 4950.75
$ws.Range("L126").Value = 7860
$ws.Range("M126").Value = -2480.75
$ws.Range("N126").Value = -12800

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1039.2727
$ws.Range("I113").Value = 1039.4166
$ws.Range("J113").Value = 1039.1
$ws.Range("K113").Value = 3118.2498
$ws.Range("L113").Value = 3117.3
$ws.Range("M113").Value = -948.2498000000001
$ws.Range("N113").Value = -7457.299999999999
